$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update "Förändrad" (C column) date for rows 2-8 from 2023-09-01 (45170) to 2023-09-05 (45174)
foreach ($r in 2..8) {
    $ws.Cells.Item($r, 3).Value = 45174
}
